# Update the demo medicine statistics table and drop the obsolete third row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: refresh stats for the first medicine ---
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = "Ab"
$ws.Range("C1").Value = "demoMed"
$ws.Range("D1").Value = 1.5
$ws.Range("E1").Value = 10
$ws.Range("F1").Value = 100
$ws.Range("G1").Value = "red"
$ws.Range("H1").Value = "shape"

# --- Row 2: refresh stats for the second medicine ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Rx"
$ws.Range("C2").Value = "Ma tuy"
$ws.Range("D2").Value = 1.2
$ws.Range("E2").Value = 25
$ws.Range("F2").Value = 165
$ws.Range("G2").Value = "Green"
$ws.Range("H2").Value = "Round"

# --- Row 3 is no longer needed, remove it entirely ---
$ws.Rows("3:3").Delete()

# Match the new active selection left behind by the edit
$null = $ws.Range("G2").Select()
